$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = 2499.5183567505
$ws.Range("I2").Value = 54.5183567504973
$ws.Range("B3").Value = 2572.81159379997
$ws.Range("I3").Value = 1639.81159379997
$ws.Range("B4").Value = 4564.30686758376
$ws.Range("C4").Value = 3516.58659342595
$ws.Range("D4").Value = 3406.01186754749
$ws.Range("I4").Value = 1056.30686758376
$ws.Range("B5").Value = 5160.7768936451
$ws.Range("C5").Value = 3754.88327627754
$ws.Range("D5").Value = 3580.64703186017
$ws.Range("I5").Value = -199.223106354897
$ws.Range("B6").Value = 6045.47843152312
$ws.Range("C6").Value = 4179.74185310249
$ws.Range("D6").Value = 3834.92948709128
$ws.Range("I6").Value = 781.478431523115
$ws.Range("B7").Value = 6400.93935985961
$ws.Range("C7").Value = 4599.86481061354
$ws.Range("D7").Value = 4163.71888182886
$ws.Range("I7").Value = 207.939359859613
$ws.Range("B8").Value = 6174.23920098852
$ws.Range("C8").Value = 4416.47475631294
$ws.Range("D8").Value = 4121.66369501346
$ws.Range("I8").Value = 1202.23920098852
$ws.Range("B9").Value = 5383.41921426064
$ws.Range("I9").Value = 411.419214260643
$ws.Range("B10").Value = 3850.32131348732
$ws.Range("E10").Value = 5368.3103318324
$ws.Range("F10").Value = 6016.43089291314
$ws.Range("I10").Value = -641.678686512675
$ws.Range("B11").Value = 2808.34359226312
$ws.Range("E11").Value = 4000.07775683539
$ws.Range("I11").Value = -397.656407736879
$ws.Range("B12").Value = 2971.73113155808
$ws.Range("I12").Value = -639.26886844192
$ws.Range("B13").Value = 3065.58137355917
$ws.Range("I13").Value = -628.418626440831
$ws.Range("B14").Value = 2600.60698409385
$ws.Range("I14").Value = -574.393015906146
$ws.Range("B15").Value = 2677.83936486815
$ws.Range("C15").Value = 1669.35358350412
$ws.Range("I15").Value = -747.160635131854
$ws.Range("B16").Value = 4412.19712534293
$ws.Range("C16").Value = 2368.88119195857
$ws.Range("D16").Value = 2144.89143594314
$ws.Range("I16").Value = -2807.80287465707
$ws.Range("B17").Value = 5188.66949872935
$ws.Range("C17").Value = 3063.85344572253
$ws.Range("D17").Value = 2705.98458571268
$ws.Range("I17").Value = -2659.33050127065
$ws.Range("B18").Value = 5934.7260396503
$ws.Range("C18").Value = 3278.04775288533
$ws.Range("D18").Value = 2811.77276626491
$ws.Range("I18").Value = -2161.2739603497
$ws.Range("B19").Value = 6241.11540961488
$ws.Range("C19").Value = 3778.3376961019
$ws.Range("D19").Value = 3353.6368979253
$ws.Range("I19").Value = -3428.88459038512
$ws.Range("B20").Value = 6125.08317192318
$ws.Range("C20").Value = 4241.43946751464
$ws.Range("I20").Value = -3096.91682807682
$ws.Range("B21").Value = 5311.8800883987
$ws.Range("I21").Value = -1555.1199116013
$ws.Range("B22").Value = 3881.57466336034
$ws.Range("E22").Value = 6300.13849596426
$ws.Range("F22").Value = 7551.56273391586
$ws.Range("I22").Value = -2050.42533663966
$ws.Range("B23").Value = 2787.78180835005
$ws.Range("E23").Value = 4866.72494197849
$ws.Range("I23").Value = -834.218191649951
$ws.Range("B24").Value = 2858.44690328899
$ws.Range("I24").Value = -790.553096711012
$ws.Range("B25").Value = 3003.72585016267
$ws.Range("C25").Value = 1679.17274067448
$ws.Range("I25").Value = -1915.27414983733
$ws.Range("B26").Value = 2578.59594913963
$ws.Range("C26").Value = 1274.39903375164
$ws.Range("I26").Value = -1817.40405086037
$ws.Range("B27").Value = 2636.29555131407
$ws.Range("C27").Value = 1246.91091622409
$ws.Range("I27").Value = -2052.70444868593
$ws.Range("B28").Value = 4246.798660274
$ws.Range("C28").Value = 1464.79754187174
$ws.Range("D28").Value = 1234.05437160573
$ws.Range("I28").Value = -2409.201339726
$ws.Range("B29").Value = 5200.73329114052
$ws.Range("C29").Value = 2142.04282848902
$ws.Range("D29").Value = 1489.17108245865
$ws.Range("I29").Value = -1923.26670885948
$ws.Range("B30").Value = 6049.68962963434
$ws.Range("C30").Value = 2998.47523940712
$ws.Range("D30").Value = 2385.62802207643
$ws.Range("I30").Value = -2774.31037036566
$ws.Range("B31").Value = 6369.63562280809
$ws.Range("C31").Value = 3677.88420378451
$ws.Range("I31").Value = -3573.36437719191
$ws.Range("B32").Value = 6289.62496318651
$ws.Range("I32").Value = -3393.37503681349
$ws.Range("B33").Value = 5583.3078550542
$ws.Range("I33").Value = -2303.6921449458
$ws.Range("B34").Value = 4026.78098672855
$ws.Range("E34").Value = 7196.2425014789
$ws.Range("I34").Value = -1284.21901327145
$ws.Range("B35").Value = 2904.29455209353
$ws.Range("E35").Value = 5670.92569293253
$ws.Range("I35").Value = 369.294552093527
$ws.Range("B36").Value = 2978.98236426928
$ws.Range("I36").Value = 409.982364269279
$ws.Range("B37").Value = 3116.12243529275
$ws.Range("C37").Value = 1376.03299995998
$ws.Range("I37").Value = 1296.12243529275
$ws.Range("B38").Value = 2707.58079303788
$ws.Range("I38").Value = 389.580793037882
$ws.Range("B39").Value = 2731.14513363903
$ws.Range("I39").Value = -2930.85486636097
$ws.Range("B40").Value = 4365.31031859338
$ws.Range("C40").Value = 1162.17051701505
$ws.Range("D40").Value = 976.220975602741
$ws.Range("I40").Value = -2177.68968140662
$ws.Range("B41").Value = 5381.17211393562
$ws.Range("C41").Value = 1471.47832539291
$ws.Range("D41").Value = 1060.14540843188
$ws.Range("I41").Value = -2295.82788606438
$ws.Range("B42").Value = 6303.01809675458
$ws.Range("C42").Value = 2882.99615025094
$ws.Range("I42").Value = -2763.98190324542
$ws.Range("B43").Value = 6612.61453180103
$ws.Range("I43").Value = -1713.38546819897
$ws.Range("B44").Value = 6566.2962390924
$ws.Range("I44").Value = -2597.7037609076
$ws.Range("B45").Value = 5794.31686845813
$ws.Range("I45").Value = -2559.68313154187
$ws.Range("B46").Value = 4210.26938612463
$ws.Range("I46").Value = -1776.73061387537
$ws.Range("B47").Value = 3083.9471661694
$ws.Range("I47").Value = -1393.0528338306
$ws.Range("B48").Value = 3113.10906076194
$ws.Range("I48").Value = -1426.89093923806
$ws.Range("B49").Value = 3251.30989715978
$ws.Range("I49").Value = -491.690102840225
